$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- In-place cell edits (rows 7, 8, 10) ---
# Row 7
$ws.Range("A7").Value = 16
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0"

# Row 8
$ws.Range("A8").Value = 16
$ws.Range("B8").Value = "Femenino"
$ws.Range("C8").Value = "No"
$ws.Range("E8").Value = 17

# Row 10
$ws.Range("A10").Value = 15
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0"

# --- Insert a new row before row 44, shifting old rows 44-75 down to 45-76 ---
$ws.Rows(44).Insert()

# The newly inserted row 44 gets the content that used to be in row 43
# (before row 43 itself is edited below).
$ws.Range("A44").Value = 13
$ws.Range("B44").Value = "Femenino"
$ws.Range("C44").Value = "Si"
$ws.Range("D44").Value = "Si"
$ws.Range("E44").Value = 16

# Row 43 is then edited to its new final values.
$ws.Range("A43").Value = 19
$ws.Range("B43").Value = "Femenino"
$ws.Range("C43").Value = "No"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1"
$ws.Range("E43").Value = 16
